# Updates the "UniformA-HW20" sheet with the new simulation-scheme data:
#  - shrinks the duplicated HKL columns (X:AG) off the header/table,
#  - relabels the HKL reflections shown in row 2,
#  - renames all the simulation schemes in column B (rows 3:19),
#  - appends 10 new simulation schemes as rows 20:29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the old duplicate columns X1:AG19 - the new table only spans A:W.
$ws.Range("X1:AG19").Clear() | Out-Null

# 2) Row 2 (HKL labels row): C2:W2 get reordered/renamed HKL + pair-scheme labels.
$row2 = New-Object 'object[,]' 1,21
$row2Vals = @('[5, 1, 1]', '[4, 2, 2]', '[3, 3, 1]', '[3, 1, 1]', '[1, 1, 1]', '[2, 2, 2]', '[3, 3, 3]', '[2, 0, 0]', '[2, 2, 0]', '[4, 2, 0]', '[4, 0, 0]', '1Pair-A', '1Pair-B', '2Pairs-A', '2Pairs-B', '3Pairs-A', '3Pairs-B', '3Pairs-C', '4Pairs', '5A4F', 'MaxUnique')
for ($i = 0; $i -lt 21; $i++) { $row2[0, $i] = $row2Vals[$i] }
$ws.Range("C2:W2").Value = $row2

# 3) Column B, rows 3:19 - existing scheme rows get renamed in place.
$b3to19 = New-Object 'object[,]' 17,1
$b3to19Vals = @('Spiral5', 'RotRing OmegaMax-90', 'Equal Angle', 'Tilt Rotate', 'CLR', 'Rizzie Hex', 'Thomas Hex', 'Tilt Rotate_Partial', 'RotRing OmegaMax-60', 'Equal Angle_Partial', 'Rizzie Hex_Partial', 'ND Single', 'RD Single', 'TD Single', 'Morris Single', 'Ring Perpendicular to ND', 'Ring Perpendicular to RD')
for ($i = 0; $i -lt 17; $i++) { $b3to19[$i, 0] = $b3to19Vals[$i] }
$ws.Range("B3:B19").Value = $b3to19

# 4) New rows 20:29 - ten additional simulation schemes, each a full 1..21 row of 1s,
#    with its own index in column A and name in column B.
$b20to29Vals = @('Ring Perpendicular to TD', 'OffsetFTD', 'OffsetATD', 'OffsetF45', 'OffsetA45', 'OffsetFRD', 'OffsetARD', 'Gaussian Quadrature', 'Michael-CCHex', 'Michael-SNHex')

$newRows = New-Object 'object[,]' 10,23
for ($i = 0; $i -lt 10; $i++) {
    $newRows[$i, 0] = 18 + $i          # column A: running index
    $newRows[$i, 1] = $b20to29Vals[$i] # column B: scheme name
    for ($j = 2; $j -lt 23; $j++) {
        $newRows[$i, $j] = 1           # columns C:W: all 1
    }
}
$ws.Range("A20:W29").Value = $newRows

# 5) Match formatting of the new cells to the existing table:
#    column A index cells + header-row style (bold, thin border, centered).
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A20:A29").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = $false

Write-Host "UniformA-HW20 sheet updated with new simulation scheme data"
